$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 43.2
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H70").Value = 2941.75
$ws.Range("I70").Value = 3143.1428
$ws.Range("J70").Value = 2659.8
$ws.Range("K70").Value = 9429.428400000001
$ws.Range("L70").Value = 7979.400000000001
$ws.Range("M70").Value = -9159.428400000001
$ws.Range("N70").Value = -8519.400000000001
$ws.Range("H73").Value = 2941.75
$ws.Range("I73").Value = 3143.1428
$ws.Range("J73").Value = 2659.8
$ws.Range("K73").Value = 9429.428400000001
$ws.Range("L73").Value = 7979.400000000001
$ws.Range("M73").Value = -8493.428400000001
$ws.Range("N73").Value = -9851.400000000001
$ws.Range("H106").Value = 51284220
$ws.Range("I106").Value = 23811732
$ws.Range("J106").Value = 83335460
$ws.Range("K106").Value = 23811732
$ws.Range("L106").Value = 83335460
$ws.Range("M106").Value = -23811101
$ws.Range("N106").Value = -83336722
$ws.Range("H132").Value = 1180.1702
$ws.Range("I132").Value = 922.9
$ws.Range("J132").Value = 2650.2856
$ws.Range("K132").Value = 2768.7
$ws.Range("L132").Value = 7950.8568
$ws.Range("M132").Value = -238.6999999999998
$ws.Range("N132").Value = -13010.8568

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 232.33333
$ws.Range("I5").Value = 260.4
$ws.Range("J5").Value = 92
$ws.Range("K5").Value = 260.4
$ws.Range("L5").Value = 92
$ws.Range("M5").Value = -148.4
$ws.Range("N5").Value = -316
$ws.Range("H32").Value = 5226.894
$ws.Range("I32").Value = 3297.8071
$ws.Range("K32").Value = 3297.8071
$ws.Range("M32").Value = -3010.8071
$ws.Range("H61").Value = 3283.68
$ws.Range("I61").Value = 3248.653
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 3248.653
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -3036.653
$ws.Range("N61").Value = -5424
$ws.Range("H74").Value = 1240.2051
$ws.Range("I74").Value = 1076.3871
$ws.Range("J74").Value = 1875
$ws.Range("K74").Value = 1076.3871
$ws.Range("L74").Value = 1875
$ws.Range("M74").Value = -202.3870999999999
$ws.Range("N74").Value = -3623
$ws.Range("H77").Value = 1240.2051
$ws.Range("I77").Value = 1076.3871
$ws.Range("J77").Value = 1875
$ws.Range("K77").Value = 5381.9355
$ws.Range("L77").Value = 9375
$ws.Range("M77").Value = -1013.9355
$ws.Range("N77").Value = -18111
$ws.Range("H110").Value = 892.1429000000001
$ws.Range("I110").Value = 892.1429000000001
$ws.Range("K110").Value = 892.1429000000001
$ws.Range("M110").Value = 1152.8571
$ws.Range("H122").Value = 2565473.8
$ws.Range("I122").Value = 2850302.8
$ws.Range("J122").Value = 2014
$ws.Range("K122").Value = 8550908.399999999
$ws.Range("L122").Value = 6042
$ws.Range("M122").Value = -8548458.399999999
$ws.Range("N122").Value = -10942
$ws.Range("H132").Value = 2055.5166
$ws.Range("I132").Value = 1179.2554
$ws.Range("J132").Value = 5223.5386
$ws.Range("K132").Value = 3537.7662
$ws.Range("L132").Value = 15670.6158
$ws.Range("M132").Value = -1007.7662
$ws.Range("N132").Value = -20730.6158
$ws.Range("H136").Value = 3283.68
$ws.Range("I136").Value = 3248.653
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 9745.958999999999
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -7195.958999999999
$ws.Range("N136").Value = -20100

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 232.33333
$ws.Range("I4").Value = 260.4
$ws.Range("J4").Value = 92
$ws.Range("K4").Value = 260.4
$ws.Range("L4").Value = 92
$ws.Range("M4").Value = -145.4
$ws.Range("N4").Value = -322
$ws.Range("H20").Value = 11164.5
$ws.Range("I20").Value = 743.86664
$ws.Range("J20").Value = 28532.223
$ws.Range("K20").Value = 743.86664
$ws.Range("L20").Value = 28532.223
$ws.Range("M20").Value = -496.86664
$ws.Range("N20").Value = -29026.223
$ws.Range("H105").Value = 1955.04
$ws.Range("I105").Value = 1540.8948
$ws.Range("J105").Value = 3266.5
$ws.Range("K105").Value = 1540.8948
$ws.Range("L105").Value = 3266.5
$ws.Range("M105").Value = 206.1052
$ws.Range("N105").Value = -6760.5
$ws.Range("H107").Value = 672
$ws.Range("I107").Value = 682.0952
$ws.Range("K107").Value = 682.0952
$ws.Range("M107").Value = 1237.9048
$ws.Range("H134").Value = 5051.8
$ws.Range("I134").Value = 5659.7407
$ws.Range("K134").Value = 16979.2221
$ws.Range("M134").Value = -14444.2221

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 120
$ws.Range("I7").Value = 71.25
$ws.Range("K7").Value = 71.25
$ws.Range("M7").Value = 41.75
$ws.Range("H31").Value = 2834.5083
$ws.Range("I31").Value = 1398.875
$ws.Range("J31").Value = 5569.048
$ws.Range("K31").Value = 1398.875
$ws.Range("L31").Value = 5569.048
$ws.Range("M31").Value = -1103.875
$ws.Range("N31").Value = -6159.048
$ws.Range("H34").Value = 2834.5083
$ws.Range("I34").Value = 1398.875
$ws.Range("J34").Value = 5569.048
$ws.Range("K34").Value = 1398.875
$ws.Range("L34").Value = 5569.048
$ws.Range("M34").Value = -1196.875
$ws.Range("N34").Value = -5973.048
$ws.Range("H62").Value = 8200
$ws.Range("I62").Value = 9266.666999999999
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 9266.666999999999
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -8642.666999999999
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 8200
$ws.Range("I65").Value = 9266.666999999999
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 46333.335
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -43213.335
$ws.Range("N65").Value = -31240

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2255.7778
$ws.Range("J22").Value = 2671.7144
$ws.Range("L22").Value = 8015.1432
$ws.Range("N22").Value = -8353.143199999999
$ws.Range("H27").Value = 2255.7778
$ws.Range("J27").Value = 2671.7144
$ws.Range("L27").Value = 8015.1432
$ws.Range("N27").Value = -8219.143199999999
$ws.Range("H44").Value = 896.11536
$ws.Range("J44").Value = 1149.9375
$ws.Range("L44").Value = 3449.8125
$ws.Range("N44").Value = -4245.8125
$ws.Range("H68").Value = 1256.6666
$ws.Range("I68").Value = 486.75
$ws.Range("J68").Value = 1730.4615
$ws.Range("K68").Value = 1460.25
$ws.Range("L68").Value = 5191.3845
$ws.Range("M68").Value = -649.25
$ws.Range("N68").Value = -6813.3845
$ws.Range("H71").Value = 1256.6666
$ws.Range("I71").Value = 486.75
$ws.Range("J71").Value = 1730.4615
$ws.Range("K71").Value = 4380.75
$ws.Range("L71").Value = 15574.1535
$ws.Range("M71").Value = -324.75
$ws.Range("N71").Value = -23686.1535
$ws.Range("H116").Value = 3000
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H123").Value = 3418.5715
$ws.Range("I123").Value = 3732.5
$ws.Range("J123").Value = 3000
$ws.Range("K123").Value = 11197.5
$ws.Range("L123").Value = 9000
$ws.Range("M123").Value = -8747.5
$ws.Range("N123").Value = -13900

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 9000
$ws.Range("J38").Value = 9000
$ws.Range("L38").Value = 9000
$ws.Range("N38").Value = -9926
$ws.Range("H102").Value = 1988.8889
$ws.Range("I102").Value = 1766.6666
$ws.Range("K102").Value = 1766.6666
$ws.Range("M102").Value = -144.6666
$ws.Range("H113").Value = 90910560
$ws.Range("I113").Value = 111112210
$ws.Range("J113").Value = 3150
$ws.Range("K113").Value = 111112210
$ws.Range("L113").Value = 3150
$ws.Range("M113").Value = -111110040
$ws.Range("N113").Value = -7490

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 2933.125
$ws.Range("I35").Value = 2277.5
$ws.Range("J35").Value = 4900
$ws.Range("K35").Value = 2277.5
$ws.Range("L35").Value = 4900
$ws.Range("M35").Value = -1941.5
$ws.Range("N35").Value = -5572
$ws.Range("H43").Value = 5000
$ws.Range("J43").Value = 5000
$ws.Range("L43").Value = 5000
$ws.Range("N43").Value = -5386
$ws.Range("H51").Value = 11032
$ws.Range("I51").Value = 1076
$ws.Range("J51").Value = 13521
$ws.Range("K51").Value = 1076
$ws.Range("L51").Value = 13521
$ws.Range("M51").Value = -598
$ws.Range("N51").Value = -14477
$ws.Range("H55").Value = 100000300
$ws.Range("I55").Value = 100000216
$ws.Range("J55").Value = 100000480
$ws.Range("K55").Value = 100000216
$ws.Range("L55").Value = 100000480
$ws.Range("M55").Value = -100000043
$ws.Range("N55").Value = -100000826
$ws.Range("H68").Value = 200003340
$ws.Range("I68").Value = 4175
$ws.Range("K68").Value = 4175
$ws.Range("M68").Value = -3426
$ws.Range("H71").Value = 200003340
$ws.Range("I71").Value = 4175
$ws.Range("K71").Value = 20875
$ws.Range("M71").Value = -17131
$ws.Range("H82").Value = 1003909.06
$ws.Range("I82").Value = 2001459.8
$ws.Range("J82").Value = 172616.83
$ws.Range("K82").Value = 2001459.8
$ws.Range("L82").Value = 172616.83
$ws.Range("M82").Value = -2001098.8
$ws.Range("N82").Value = -173338.83
$ws.Range("H85").Value = 1003909.06
$ws.Range("I85").Value = 2001459.8
$ws.Range("J85").Value = 172616.83
$ws.Range("K85").Value = 2001459.8
$ws.Range("L85").Value = 172616.83
$ws.Range("M85").Value = -2000211.8
$ws.Range("N85").Value = -175112.83
$ws.Range("H133").Value = 39791.25
$ws.Range("J133").Value = 39791.25
$ws.Range("L133").Value = 39791.25
$ws.Range("N133").Value = -44851.25

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1512.1111
$ws.Range("I122").Value = 1413
$ws.Range("J122").Value = 2305
$ws.Range("K122").Value = 4239
$ws.Range("L122").Value = 6915
$ws.Range("M122").Value = -1789
$ws.Range("N122").Value = -11815
$ws.Range("H136").Value = 1165.4524
$ws.Range("I136").Value = 682.7778
$ws.Range("J136").Value = 1527.4584
$ws.Range("K136").Value = 2048.3334
$ws.Range("L136").Value = 4582.3752
$ws.Range("M136").Value = 501.6666
$ws.Range("N136").Value = -9682.3752
